{"js": "const body = context.document.body;\n\n// ---- Change 1: merge \"Within the game, rules of play will be available.\" into a single run ----\n{\n  const results = body.search(\"Within the game, rules of play will be available.\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Within the game, rules of play will be available.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---- Change 2: \"...and start the timer clock.\" -> \"...and display the timer clock.\" (REQ-2, section 3.2.2) ----\n{\n  const results = body.search(\"and start the timer clock\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"and display the timer clock\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- Change 3: \"...and record the time from the clock.\" -> \"...and display the ending time from the clock.\" (REQ-4, section 3.2.2) ----\n{\n  const results = body.search(\"and record the time from the clock\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"and display the ending time from the clock\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---- Change 4: remove the \"high scores\" REQ-2 paragraph (section 3.4.2) and renumber the following two requirements ----\n{\n  const results = body.search(\n    \"The system should keep track of the high scores for each solitaire variant.\",\n    { matchCase: true }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const paraCollection = results.items[0].paragraphs;\n    paraCollection.load(\"items\");\n    await context.sync();\n    paraCollection.items[0].delete();\n    await context.sync();\n  }\n}\n\n// Old \"REQ-3:\" (display current score/time) becomes \"REQ-2:\"\n{\n  const results = body.search(\n    \"The system should display the current score and time while the game is being played.\",\n    { matchCase: true }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const paraCollection = results.items[0].paragraphs;\n    paraCollection.load(\"items\");\n    await context.sync();\n    const labelRange = paraCollection.items[0].getRange();\n    const labelResults = labelRange.search(\"REQ-3:\", { matchCase: true });\n    labelResults.load(\"text\");\n    await context.sync();\n    if (labelResults.items.length > 0) {\n      labelResults.items[0].insertText(\"REQ-2:\", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n\n// Old \"REQ-4:\" (keep track of time spent until finish) becomes \"REQ-3:\"\n{\n  const results = body.search(\n    \"The system should keep track of the current time spent from the start of the game\",\n    { matchCase: true }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const paraCollection = results.items[0].paragraphs;\n    paraCollection.load(\"items\");\n    await context.sync();\n    const labelRange = paraCollection.items[0].getRange();\n    const labelResults = labelRange.search(\"REQ-4:\", { matchCase: true });\n    labelResults.load(\"text\");\n    await context.sync();\n    if (labelResults.items.length > 0) {\n      labelResults.items[0].insertText(\"REQ-3:\", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---- Change 1: merge \"Within the game, rules of play will be available.\" into a single run ----\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Replacement.ClearFormatting()\n$r1.Find.Text = \"Within the game, rules of play will be available.\"\n$r1.Find.Replacement.Text = \"Within the game, rules of play will be available.\"\n$r1.Find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n\n# ---- Change 2: \"...and start the timer clock.\" -> \"...and display the timer clock.\" (REQ-2, section 3.2.2) ----\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Replacement.ClearFormatting()\n$r2.Find.Text = \"and start the timer clock\"\n$r2.Find.Replacement.Text = \"and display the timer clock\"\n$r2.Find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n\n# ---- Change 3: \"...and record the time from the clock.\" -> \"...and display the ending time from the clock.\" (REQ-4, section 3.2.2) ----\n$r3 = $d.Content\n$r3.Find.ClearFormatting()\n$r3.Find.Replacement.ClearFormatting()\n$r3.Find.Text = \"and record the time from the clock\"\n$r3.Find.Replacement.Text = \"and display the ending time from the clock\"\n$r3.Find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n\n# ---- Change 4: remove the \"high scores\" REQ-2 paragraph (section 3.4.2) and renumber the following two requirements ----\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*high scores for each solitaire variant*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n\n    # Old \"REQ-3:\" (display current score/time) becomes \"REQ-2:\"\n    $pA = $d.Paragraphs.Item($targetIndex)\n    $rA = $pA.Range\n    $rA.Find.ClearFormatting()\n    $rA.Find.Replacement.ClearFormatting()\n    $rA.Find.Text = \"REQ-3:\"\n    $rA.Find.Replacement.Text = \"REQ-2:\"\n    $rA.Find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n\n    # Old \"REQ-4:\" (keep track of time spent until finish) becomes \"REQ-3:\"\n    $pB = $d.Paragraphs.Item($targetIndex + 1)\n    $rB = $pB.Range\n    $rB.Find.ClearFormatting()\n    $rB.Find.Replacement.ClearFormatting()\n    $rB.Find.Text = \"REQ-4:\"\n    $rB.Find.Replacement.Text = \"REQ-3:\"\n    $rB.Find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n}\n"}
